$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 57, shifting existing rows 57-68 down to 58-69.
$ws.Rows.Item(57).Insert()

# Populate the newly inserted row 57 with the new weekly record.
$ws.Cells.Item(57, 1).Value = 4
$ws.Cells.Item(57, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(57, 3).Value = "Los Lagos"
$ws.Cells.Item(57, 4).NumberFormat = $ws.Cells.Item(58, 4).NumberFormat
$ws.Cells.Item(57, 4).Value = 45093
$ws.Cells.Item(57, 5).Value = 10
$ws.Cells.Item(57, 6).Value = 100112012
$ws.Cells.Item(57, 7).Value = "Espinaca"
$ws.Cells.Item(57, 8).Value = "Sin especificar"
$ws.Cells.Item(57, 9).Value = "Primera"
$ws.Cells.Item(57, 10).Value = 25
$ws.Cells.Item(57, 11).Value = 13000
$ws.Cells.Item(57, 12).Value = 13000
$ws.Cells.Item(57, 13).Value = 13000
$ws.Cells.Item(57, 14).Value = "`$/cuna 10 kilos"
$ws.Cells.Item(57, 15).Value = "Región Metropolitana"
$ws.Cells.Item(57, 16).Value = 1300
$ws.Cells.Item(57, 17).Value = 10
$ws.Cells.Item(57, 18).Value = "Hortaliza"
